# Fruta / hortaliza, semanal
# Insert 4 new weekly records (Bing/Santina, Primera/Segunda) for
# Vega Monumental Concepción - Cereza, pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows above the current row 56 - this shifts the existing
# rows 56:120 down to 60:124 (formatting of row 56/D column is carried
# down from the row above, matching the original style for the date column).
$ws.Rows("56:59").Insert()

# New row 56: Bing, Primera
$ws.Cells.Item(56, 1).Value = 11
$ws.Cells.Item(56, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(56, 3).Value = "Bíobío"
$ws.Cells.Item(56, 4).Value = 44894
$ws.Cells.Item(56, 5).Value = 8
$ws.Cells.Item(56, 6).Value = "Fruta"
$ws.Cells.Item(56, 7).Value = 100103
$ws.Cells.Item(56, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(56, 9).Value = 100103001
$ws.Cells.Item(56, 10).Value = "Cereza"
$ws.Cells.Item(56, 11).Value = "Bing"
$ws.Cells.Item(56, 12).Value = "Primera"
$ws.Cells.Item(56, 13).Value = 100
$ws.Cells.Item(56, 14).Value = 9000
$ws.Cells.Item(56, 15).Value = 10000
$ws.Cells.Item(56, 16).Value = 9500
$ws.Cells.Item(56, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(56, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(56, 19).Value = 950
$ws.Cells.Item(56, 20).Value = 10

# New row 57: Bing, Segunda
$ws.Cells.Item(57, 1).Value = 11
$ws.Cells.Item(57, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(57, 3).Value = "Bíobío"
$ws.Cells.Item(57, 4).Value = 44894
$ws.Cells.Item(57, 5).Value = 8
$ws.Cells.Item(57, 6).Value = "Fruta"
$ws.Cells.Item(57, 7).Value = 100103
$ws.Cells.Item(57, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(57, 9).Value = 100103001
$ws.Cells.Item(57, 10).Value = "Cereza"
$ws.Cells.Item(57, 11).Value = "Bing"
$ws.Cells.Item(57, 12).Value = "Segunda"
$ws.Cells.Item(57, 13).Value = 50
$ws.Cells.Item(57, 14).Value = 8000
$ws.Cells.Item(57, 15).Value = 8000
$ws.Cells.Item(57, 16).Value = 8000
$ws.Cells.Item(57, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(57, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(57, 19).Value = 800
$ws.Cells.Item(57, 20).Value = 10

# New row 58: Santina, Primera
$ws.Cells.Item(58, 1).Value = 11
$ws.Cells.Item(58, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(58, 3).Value = "Bíobío"
$ws.Cells.Item(58, 4).Value = 44894
$ws.Cells.Item(58, 5).Value = 8
$ws.Cells.Item(58, 6).Value = "Fruta"
$ws.Cells.Item(58, 7).Value = 100103
$ws.Cells.Item(58, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(58, 9).Value = 100103001
$ws.Cells.Item(58, 10).Value = "Cereza"
$ws.Cells.Item(58, 11).Value = "Santina"
$ws.Cells.Item(58, 12).Value = "Primera"
$ws.Cells.Item(58, 13).Value = 200
$ws.Cells.Item(58, 14).Value = 9000
$ws.Cells.Item(58, 15).Value = 10000
$ws.Cells.Item(58, 16).Value = 9500
$ws.Cells.Item(58, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(58, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(58, 19).Value = 950
$ws.Cells.Item(58, 20).Value = 10

# New row 59: Santina, Segunda
$ws.Cells.Item(59, 1).Value = 11
$ws.Cells.Item(59, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(59, 3).Value = "Bíobío"
$ws.Cells.Item(59, 4).Value = 44894
$ws.Cells.Item(59, 5).Value = 8
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100103
$ws.Cells.Item(59, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(59, 9).Value = 100103001
$ws.Cells.Item(59, 10).Value = "Cereza"
$ws.Cells.Item(59, 11).Value = "Santina"
$ws.Cells.Item(59, 12).Value = "Segunda"
$ws.Cells.Item(59, 13).Value = 100
$ws.Cells.Item(59, 14).Value = 8000
$ws.Cells.Item(59, 15).Value = 8000
$ws.Cells.Item(59, 16).Value = 8000
$ws.Cells.Item(59, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(59, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(59, 19).Value = 800
$ws.Cells.Item(59, 20).Value = 10

# Make sure the date cells use the same date number format as the rest of
# column D (matches the "s=2" style already used by every other row).
$ws.Range("D56:D59").NumberFormat = "YYYY-MM-DD HH:MM:SS"
